# edit.ps1
# Applies the "sent to P for forecasts" changes:
#  - renames "Input - forecasts" -> "Input - forecast"
#  - adds 4 quarters (2017Q1..2017Q4) of forecast input data
#  - extends "Escsount" sheet with matching escalation rows
#  - extends "Output" sheet with matching forecast rows (replacing the old
#    placeholder rows that simply repeated the last history row)
#  - restores a few view/selection states

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the "Input - forecasts" sheet to "Input - forecast"
# ---------------------------------------------------------------------------
$wsFore = $wb.Worksheets.Item("Input - forecasts")
$wsFore.Name = "Input - forecast"

# ---------------------------------------------------------------------------
# 2. "Input - forecast" sheet: drop the old HPIc header cell and populate
#    four quarters of forecast data (rows 2-5)
# ---------------------------------------------------------------------------
$wsFore.Range("H1").Clear()

$foreRows = @(
    @("2017Q1", 242.821, 16889.1, 16266.8, 4.7, 104.438, 185.54, 0.54),
    @("2017Q2", 242.821, 16889.1, 16266.8, 4.7, 104.438, 185.54, 0.54),
    @("2017Q3", 242.821, 16889.1, 16266.8, 4.7, 104.438, 185.54, 0.54),
    @("2017Q4", 242.821, 16889.1, 16266.8, 4.7, 104.438, 185.54, 0.54)
)

for ($i = 0; $i -lt $foreRows.Count; $i++) {
    $r = 2 + $i
    $row = $foreRows[$i]
    $wsFore.Cells.Item($r, 1).Value = $row[0]
    $wsFore.Cells.Item($r, 2).Value = $row[1]
    $wsFore.Cells.Item($r, 3).Value = $row[2]
    $wsFore.Cells.Item($r, 4).Value = $row[3]
    $wsFore.Cells.Item($r, 5).Value = $row[4]
    $wsFore.Cells.Item($r, 6).Value = $row[5]
    $wsFore.Cells.Item($r, 7).Value = $row[6]
    $wsFore.Cells.Item($r, 9).Value = $row[7]
}

$wsFore.Range("A5").Select()

# ---------------------------------------------------------------------------
# 3. "Escsount" sheet: add rows 106-109 referencing the new forecast rows
# ---------------------------------------------------------------------------
$wsEsc = $wb.Worksheets.Item("Escsount")

for ($i = 0; $i -lt 4; $i++) {
    $escRow = 106 + $i
    $foreRow = 2 + $i
    $wsEsc.Cells.Item($escRow, 1).Formula = "='Input - forecast'!A$foreRow"
    $wsEsc.Cells.Item($escRow, 2).Formula = "='Input - history'!`$B`$105/'Input - forecast'!B$foreRow"
}

$wsEsc.Range("A107").Select()

# ---------------------------------------------------------------------------
# 4. "Output" sheet: replace the placeholder rows 106-109 with formulas that
#    pull from the new "Input - forecast" rows (2-5), mirroring the pattern
#    used for the historical rows above them.
# ---------------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("Output")

for ($i = 0; $i -lt 4; $i++) {
    $outRow = 106 + $i
    $foreRow = 2 + $i

    $wsOut.Cells.Item($outRow, 2).Formula = "=IF('Input - forecast'!C$foreRow=`"`",`"`",LN('Input - forecast'!C$foreRow))"
    $wsOut.Cells.Item($outRow, 3).Formula = "=IF('Input - forecast'!D$foreRow=`"`",`"`",LN('Input - forecast'!D$foreRow))"
    $wsOut.Cells.Item($outRow, 4).Formula = "=IF('Input - forecast'!E$foreRow=`"`",`"`",LN('Input - forecast'!E$foreRow))"
    $wsOut.Cells.Item($outRow, 5).Formula = "=IF('Input - forecast'!F$foreRow=`"`",`"`",LN('Input - forecast'!F$foreRow*Escsount!`$B$outRow))"
    $wsOut.Cells.Item($outRow, 6).Formula = "=IF('Input - forecast'!G$foreRow=`"`",`"`",LN('Input - forecast'!G$foreRow))"
    $wsOut.Cells.Item($outRow, 7).Formula = "=IF('Input - forecast'!H$foreRow=`"`",`"`",LN('Input - forecast'!H$foreRow*Escsount!`$B$outRow))"
    $wsOut.Cells.Item($outRow, 8).Formula = "=IF('Input - forecast'!I$foreRow=`"`",`"`",'Input - forecast'!I$foreRow)"

    if ($outRow -eq 106) {
        $wsOut.Range("I106:L106").Clear()
    } else {
        $wsOut.Range("I$outRow`:L$outRow").ClearContents()
    }
}

$wsOut.Activate()
$wsOut.Range("B107").Select()

# ---------------------------------------------------------------------------
# 5. Misc view/selection restoration
# ---------------------------------------------------------------------------
$wsHist = $wb.Worksheets.Item("Input - history")
$wsHist.Range("B1").Select()

$wsOut.Activate()
